$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 258.18182
$ws.Range("I33").Value = 199.52632
$ws.Range("K33").Value = 199.52632
$ws.Range("M33").Value = 29.47368
$ws.Range("H41").Value = 111841.445
$ws.Range("I41").Value = 824.4286
$ws.Range("J41").Value = 500401
$ws.Range("K41").Value = 824.4286
$ws.Range("L41").Value = 500401
$ws.Range("M41").Value = -384.4286
$ws.Range("N41").Value = -501281
$ws.Range("H114").Value = 78518.336
$ws.Range("J114").Value = 78518.336
$ws.Range("L114").Value = 78518.336
$ws.Range("N114").Value = -87196.336
$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 60000
$ws.Range("L124").Value = 60000
$ws.Range("N124").Value = -69820
$ws.Range("H126").Value = 80000
$ws.Range("J126").Value = 80000
$ws.Range("L126").Value = 80000
$ws.Range("N126").Value = -89880
$ws.Range("H130").Value = 80250
$ws.Range("I130").Value = 60000
$ws.Range("J130").Value = 87000
$ws.Range("K130").Value = 60000
$ws.Range("L130").Value = 87000
$ws.Range("M130").Value = -54980
$ws.Range("N130").Value = -97040
$ws.Range("H132").Value = 2610.932
$ws.Range("I132").Value = 839.4211
$ws.Range("K132").Value = 2518.2633
$ws.Range("M132").Value = 11.73669999999993
$ws.Range("H135").Value = 4507
$ws.Range("I135").Value = 4507
$ws.Range("K135").Value = 40563
$ws.Range("M135").Value = -38028
$ws.Range("H137").Value = 1998
$ws.Range("I137").Value = 1096.4584
$ws.Range("J137").Value = 3965
$ws.Range("K137").Value = 3289.3752
$ws.Range("L137").Value = 11895
$ws.Range("M137").Value = -739.3751999999999
$ws.Range("N137").Value = -16995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1807.9231
$ws.Range("I61").Value = 1499.6
$ws.Range("K61").Value = 1499.6
$ws.Range("M61").Value = -1287.6
$ws.Range("H74").Value = 2515.75
$ws.Range("I74").Value = 1499.5264
$ws.Range("K74").Value = 1499.5264
$ws.Range("M74").Value = -625.5264
$ws.Range("H77").Value = 2515.75
$ws.Range("I77").Value = 1499.5264
$ws.Range("K77").Value = 7497.632
$ws.Range("M77").Value = -3129.632
$ws.Range("H132").Value = 2513.361
$ws.Range("I132").Value = 2374.1875
$ws.Range("K132").Value = 7122.5625
$ws.Range("M132").Value = -4592.5625
$ws.Range("H136").Value = 1807.9231
$ws.Range("I136").Value = 1499.6
$ws.Range("K136").Value = 4498.799999999999
$ws.Range("M136").Value = -1948.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 22599.8
$ws.Range("I26").Value = 16999.75
$ws.Range("K26").Value = 16999.75
$ws.Range("M26").Value = -16707.75
$ws.Range("H80").Value = 1992.2307
$ws.Range("I80").Value = 372.33334
$ws.Range("J80").Value = 2478.2
$ws.Range("K80").Value = 372.33334
$ws.Range("L80").Value = 2478.2
$ws.Range("M80").Value = 625.66666
$ws.Range("N80").Value = -4474.2
$ws.Range("H83").Value = 1992.2307
$ws.Range("I83").Value = 372.33334
$ws.Range("J83").Value = 2478.2
$ws.Range("K83").Value = 1861.6667
$ws.Range("L83").Value = 12391
$ws.Range("M83").Value = 3130.3333
$ws.Range("N83").Value = -22375
$ws.Range("H105").Value = 1438.3334
$ws.Range("I105").Value = 1483.0769
$ws.Range("K105").Value = 1483.0769
$ws.Range("M105").Value = 263.9231
$ws.Range("H134").Value = 3167.9434
$ws.Range("I134").Value = 2674.319
$ws.Range("J134").Value = 7034.6665
$ws.Range("K134").Value = 8022.957
$ws.Range("L134").Value = 21103.9995
$ws.Range("M134").Value = -5487.957
$ws.Range("N134").Value = -26173.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1083.25
$ws.Range("I16").Value = 988
$ws.Range("J16").Value = 1750
$ws.Range("K16").Value = 988
$ws.Range("L16").Value = 1750
$ws.Range("M16").Value = -701
$ws.Range("N16").Value = -2324
$ws.Range("H22").Value = 1352.4117
$ws.Range("J22").Value = 1562.125
$ws.Range("L22").Value = 1562.125
$ws.Range("N22").Value = -2262.125
$ws.Range("H31").Value = 6443.543
$ws.Range("I31").Value = 5601.4346
$ws.Range("K31").Value = 5601.4346
$ws.Range("M31").Value = -5306.4346
$ws.Range("H34").Value = 6443.543
$ws.Range("I34").Value = 5601.4346
$ws.Range("K34").Value = 5601.4346
$ws.Range("M34").Value = -5399.4346
$ws.Range("H113").Value = 1083.25
$ws.Range("I113").Value = 988
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 988
$ws.Range("L113").Value = 1750
$ws.Range("M113").Value = 1182
$ws.Range("N113").Value = -6090
$ws.Range("H122").Value = 4438.875
$ws.Range("I122").Value = 4099.4
$ws.Range("J122").Value = 5004.6665
$ws.Range("K122").Value = 12298.2
$ws.Range("L122").Value = 15013.9995
$ws.Range("M122").Value = -9848.199999999999
$ws.Range("N122").Value = -19913.9995
$ws.Range("H125").Value = 65000
$ws.Range("J125").Value = 65000
$ws.Range("L125").Value = 65000
$ws.Range("N125").Value = -69920
$ws.Range("H132").Value = 173271.17
$ws.Range("I132").Value = 809.7778
$ws.Range("K132").Value = 2429.3334
$ws.Range("M132").Value = 100.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 27.0625
$ws.Range("I12").Value = 66.666664
$ws.Range("J12").Value = 17.923077
$ws.Range("K12").Value = 199.999992
$ws.Range("L12").Value = 53.769231
$ws.Range("M12").Value = -26.99999199999999
$ws.Range("N12").Value = -399.769231
$ws.Range("H38").Value = 15.4
$ws.Range("I38").Value = 5.3333335
$ws.Range("K38").Value = 16.0000005
$ws.Range("M38").Value = 330.9999995
$ws.Range("H70").Value = 790.55554
$ws.Range("I70").Value = 395.85715
$ws.Range("J70").Value = 2172
$ws.Range("K70").Value = 1187.57145
$ws.Range("L70").Value = 6516
$ws.Range("M70").Value = -872.5714499999999
$ws.Range("N70").Value = -7146
$ws.Range("H73").Value = 790.55554
$ws.Range("I73").Value = 395.85715
$ws.Range("J73").Value = 2172
$ws.Range("K73").Value = 1187.57145
$ws.Range("L73").Value = 6516
$ws.Range("M73").Value = -95.57144999999991
$ws.Range("N73").Value = -8700
$ws.Range("H75").Value = 5116.75
$ws.Range("I75").Value = 2253.75
$ws.Range("J75").Value = 6071.0835
$ws.Range("K75").Value = 6761.25
$ws.Range("L75").Value = 18213.2505
$ws.Range("M75").Value = -5763.25
$ws.Range("N75").Value = -20209.2505
$ws.Range("H78").Value = 5116.75
$ws.Range("I78").Value = 2253.75
$ws.Range("J78").Value = 6071.0835
$ws.Range("K78").Value = 20283.75
$ws.Range("L78").Value = 54639.7515
$ws.Range("M78").Value = -15291.75
$ws.Range("N78").Value = -64623.7515

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 288.46875
$ws.Range("I2").Value = 199.63158
$ws.Range("J2").Value = 418.30768
$ws.Range("K2").Value = 199.63158
$ws.Range("L2").Value = 418.30768
$ws.Range("M2").Value = -86.63158000000001
$ws.Range("N2").Value = -644.30768
$ws.Range("H102").Value = 2114.7358
$ws.Range("I102").Value = 2137.9333
$ws.Range("K102").Value = 2137.9333
$ws.Range("M102").Value = -515.9333000000001
$ws.Range("H126").Value = 2955.5557
$ws.Range("I126").Value = 2842.8572
$ws.Range("K126").Value = 8528.571599999999
$ws.Range("M126").Value = -6058.571599999999
$ws.Range("H132").Value = 23266798
$ws.Range("I132").Value = 25008010
$ws.Range("K132").Value = 75024030
$ws.Range("M132").Value = -75021500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1299.6666
$ws.Range("I22").Value = 798
$ws.Range("K22").Value = 798
$ws.Range("M22").Value = -503
$ws.Range("H27").Value = 1299.6666
$ws.Range("I27").Value = 798
$ws.Range("K27").Value = 798
$ws.Range("M27").Value = -691
$ws.Range("H41").Value = 35000
$ws.Range("I41").Value = 35000
$ws.Range("K41").Value = 35000
$ws.Range("M41").Value = -34562
$ws.Range("H46").Value = 1535.3572
$ws.Range("I46").Value = 1215.8334
$ws.Range("J46").Value = 1775
$ws.Range("K46").Value = 1215.8334
$ws.Range("L46").Value = 1775
$ws.Range("M46").Value = -1027.8334
$ws.Range("N46").Value = -2151
$ws.Range("H101").Value = 63208.168
$ws.Range("J101").Value = 63208.168
$ws.Range("L101").Value = 63208.168
$ws.Range("N101").Value = -69698.16800000001
$ws.Range("H132").Value = 2370.5715
$ws.Range("I132").Value = 2308.9614
$ws.Range("J132").Value = 3171.5
$ws.Range("K132").Value = 6926.8842
$ws.Range("L132").Value = 9514.5
$ws.Range("M132").Value = -4396.8842
$ws.Range("N132").Value = -14574.5
$ws.Range("H136").Value = 2876.7778
$ws.Range("I136").Value = 1627.4286
$ws.Range("K136").Value = 4882.2858
$ws.Range("M136").Value = -2332.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 59807
$ws.Range("J46").Value = 59807
$ws.Range("L46").Value = 59807
$ws.Range("N46").Value = -60269
$ws.Range("H132").Value = 2240.0688
$ws.Range("I132").Value = 1734.36
$ws.Range("K132").Value = 5203.08
$ws.Range("M132").Value = -2673.08
$ws.Range("H134").Value = 59807
$ws.Range("J134").Value = 59807
$ws.Range("L134").Value = 179421
$ws.Range("N134").Value = -184491
$ws.Range("H136").Value = 2873.6606
$ws.Range("I136").Value = 2384.745
$ws.Range("K136").Value = 7154.235
$ws.Range("M136").Value = -4604.235
